$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("DatosGenerales")
$wsPre   = $wb.Worksheets.Item("Precondiciones")
$wsPasos = $wb.Worksheets.Item("Pasos")
$wsCam   = $wb.Worksheets.Item("Control de cambios")

# --- Precondiciones: update the precondition texts for the new "Cordoba" scenario ---
$wsPre.Range("B2").Value = '"Córdoba" Es una ciudad cargada en la base de datos, y tiene playas de estacionamientos asociadas. Algunas de las cuales aceptan tipo de vehiculo <CPA_TipoVehiculo1>, el precio para <CPA_TipoVehiculo1> esta en el rango 0 - 99, atienden los dias <CPA_DiasDeAtencion1>, en el horario 00:00 - 23:59'
$wsPre.Range("B3").Value = '<CPA_TipoPlaya1> Es un tipo de playa cargado en la base de datos'
$wsPre.Range("B4").Value = '<CPA_TipoVehiculo1> Es un tipo de vehiculo cargado en la base de datos'
$wsPre.Range("B5").Value = ' <CPA_DiasDeAtencion1>  Es un dia de atencion cargado en la base de datos.'

# --- Pasos: update the step texts/expected values with concrete data ---
$wsPasos.Range("B3").Value  = 'Ingreso "Córdoba" en el campo nombre de ciudad'
$wsPasos.Range("C4").Value  = 'Se carga la pagina BuscarPlayas, con todas las playas de "Córdoba" disponibles en un mapa. Se carga la informacion de las playas en la grilla de playas debajo del mapa.'
$wsPasos.Range("B5").Value  = 'Selecciono <CPA_TipoPlaya1> en el campo Tipo de Playa'
$wsPasos.Range("B6").Value  = 'Selecciono <CPA_TipoVehiculo1> en el campo Tipo de Vehiculo'
$wsPasos.Range("B7").Value  = 'Ingreso 0 en el campo precio desde'
$wsPasos.Range("B8").Value  = 'Ingreso 99 en el campo precio hasta'
$wsPasos.Range("B9").Value  = 'Ingreso 00:00 en el campo hora desde'
$wsPasos.Range("B10").Value = 'Ingreso 23:59 en el campo hora hasta'
$wsPasos.Range("B11").Value = 'Selecciono <CPA_DiasDeAtencion1> en el campo Dias de atencion'
$wsPasos.Range("C12").Value = 'Se muestran en el mapa las playas que cumplen con los filtros seleccionados'

# Row 12 wraps less text now (shorter expected result), so it fits in two lines instead of three.
$wsPasos.Rows.Item(12).RowHeight = 25.5

# --- Control de cambios: responsible person text ---
$wsCam.Range("D2").Value = "Leonel Romero [autor]"

# --- Selections: move the active cell on each sheet without changing which tab is active ---
$wsDatos.Range("B4").Select()
$wsPre.Range("B5").Select()
$wsPasos.Range("C12").Select()

# Restore the originally active sheet/tab (Pasos) and its selection.
$wsPasos.Activate()
$wsPasos.Range("C12").Select()
